$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2"=20.32606987536358; "D2"=4.131994058382846; "E2"=10.30731302176281; "F2"=54.05460082667306; "G2"=3.784140366967575; "J2"=10.02047342394592; "K2"=19.30439556540139; "L2"=12.46610746182074; "N2"=25.47340129136326
    "B3"=20.26580694792097; "D3"=4.111851467595956; "E3"=10.30045843268443; "F3"=53.99040728273752; "G3"=3.787751041860259; "J3"=10.02586367475392; "K3"=19.16802028621895; "L3"=12.45793808157025; "N3"=25.51428359063351
    "B4"=20.23319296933037; "D4"=4.099252662659322; "E4"=10.29625916694828; "F4"=53.96112625322006; "G4"=3.790084058908104; "J4"=10.02938897333782; "K4"=19.08899155541685; "L4"=12.45539923531187; "N4"=25.54127003434116
    "B5"=20.22101534782628; "D5"=4.094060524727016; "E5"=10.29455004093529; "F5"=53.95174658536413; "G5"=3.791064070247564; "J5"=10.03087988209736; "K5"=19.05799935401901; "L5"=12.45498878786305; "N5"=25.55274119188558
    "B6"=20.21906073097486; "D6"=4.093194872674814; "E6"=10.29426636123088; "F6"=53.95034336130434; "G6"=3.791228572439502; "J6"=10.03113072991196; "K6"=19.05292717373231; "L6"=12.45495836003428; "N6"=25.55467459748633
    "B7"=20.23302422015778; "D7"=4.099182874439028; "E7"=10.29623610887071; "F7"=53.96098941615995; "G7"=3.790097156957723; "J7"=10.02940886018101; "K7"=19.08856863688198; "L7"=12.4553911713443; "N7"=25.54142281928184
    "B8"=20.30438660466874; "D8"=4.125096937129622; "E8"=10.30494742457775; "F8"=54.03036662240316; "G8"=3.785361306248874; "J8"=10.02228727608125; "K8"=19.25641478973798; "L8"=12.46277736831656; "N8"=25.4871064164144
    "B9"=20.47867353263175; "D9"=4.174076531206211; "E9"=10.32211862672088; "F9"=54.24657182407713; "G9"=3.776990241660795; "J9"=10.01002898553656; "K9"=19.62151637131672; "L9"=12.49684865131933; "N9"=25.39554183623657
    "B10"=20.62696017196143; "D10"=4.208938574985711; "E10"=10.33480737324965; "F10"=54.45385197475675; "G10"=3.771391588176609; "J10"=10.00205792493336; "K10"=19.90968819892825; "L10"=12.53370978561436; "N10"=25.3373788389627
    "B11"=20.69864258427398; "D11"=4.224553637214642; "E11"=10.3405998901691; "F11"=54.55853917058197; "G11"=3.768962931846018; "J11"=9.998655286727466; "K11"=20.04465662618948; "L11"=12.55301607922265; "N11"=25.31289640458725
    "B12"=20.7263784602863; "D12"=4.230431450424604; "E12"=10.34279661040795; "F12"=54.59966237346649; "G12"=3.768060148548126; "J12"=9.997398841868508; "K12"=20.09628307812969; "L12"=12.56068852630348; "N12"=25.303909733781
    "B13"=20.72037901585159; "D13"=4.229167130267002; "E13"=10.34232335949647; "F13"=54.59074015846129; "G13"=3.768253829222424; "J13"=9.997668014945392; "K13"=20.0851420696948; "L13"=12.55902010430947; "N13"=25.30583252794763
    "B14"=20.70091267817399; "D14"=4.225037916547446; "E14"=10.34078054944005; "F14"=54.56189281891407; "G14"=3.768888321272426; "J14"=9.99855127610626; "K14"=20.04889384384308; "L14"=12.55364007030356; "N14"=25.312151367336
    "B15"=20.68906548650564; "D15"=4.222504054250722; "E15"=10.33983596182487; "F15"=54.54441533413553; "G15"=3.769279163368618; "J15"=9.999096472247686; "K15"=20.0267568415623; "L15"=12.55039162481693; "N15"=25.31605886636602
    "B16"=20.62235921664992; "D16"=4.207913169506432; "E16"=10.33442925312876; "F16"=54.4472184015627; "G16"=3.771552676424613; "J16"=10.00228478547511; "K16"=19.90094225930619; "L16"=12.53249886355166; "N16"=25.33901859989234
    "B17"=20.58250811252914; "D17"=4.198899532755964; "E17"=10.33111783169201; "F17"=54.39024387985884; "G17"=3.772977604842196; "J17"=10.00429788985059; "K17"=19.82472354252711; "L17"=12.52216997972691; "N17"=25.35360988764991
    "B18"=20.55998563374177; "D18"=4.193692187168479; "E18"=10.3292150866671; "F18"=54.35845297721345; "G18"=3.773808317443492; "J18"=10.00547680925713; "K18"=19.78125203930118; "L18"=12.51646819178742; "N18"=25.36218841233744
    "B19"=20.55242889169779; "D19"=4.1919251356059; "E19"=10.32857116677332; "F19"=54.34785773231353; "G19"=3.7740914973737; "J19"=10.00587958604966; "K19"=19.7665976778906; "L19"=12.51457882930859; "N19"=25.36512489590278
    "B20"=20.58670916991979; "D20"=4.199861425407414; "E20"=10.33147013965034; "F20"=54.39620765579497; "G20"=3.772824767450043; "J20"=10.00408141491059; "K20"=19.83279941401113; "L20"=12.52324478317895; "N20"=25.35203736896241
    "B21"=20.70661450098507; "D21"=4.22625172565603; "E21"=10.34123362041219; "F21"=54.57032592539043; "G21"=3.768701497820908; "J21"=9.998290971328535; "K21"=20.0595271287892; "L21"=12.55521053252579; "N21"=25.31028765431299
    "B22"=20.78841787716218; "D22"=4.243293620641621; "E22"=10.34763352501776; "F22"=54.69274357311475; "G22"=3.76610514003508; "J22"=9.994693423261925; "K22"=20.21069985311995; "L22"=12.57820769216909; "N22"=25.28465895967675
    "B23"=20.74444892609664; "D23"=4.234216930990808; "E23"=10.34421594780379; "F23"=54.62662327237884; "G23"=3.76748189118655; "J23"=9.996596428844516; "K23"=20.12975593971177; "L23"=12.5657422169369; "N23"=25.29818580355938
    "B24"=20.58480866065355; "D24"=4.199426632511508; "E24"=10.33131085804491; "F24"=54.39350842721959; "G24"=3.772893829398431; "J24"=10.00417921607988; "K24"=19.82914722652108; "L24"=12.52275812784928; "N24"=25.35274771342482
    "B25"=20.42791524179625; "D25"=4.161023256871555; "E25"=10.31746192418786; "F25"=54.17953891236896; "G25"=3.779157485780831; "J25"=10.01316300516072; "K25"=19.51909588380896; "L25"=12.48554460780433; "N25"=25.41871249580331
}

foreach ($cellRef in $data.Keys) {
    $ws.Range($cellRef).Value = $data[$cellRef]
}

Write-Output "Done updating $($data.Count) cells"